$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2290180
$ws.Range("J40").Value = 3954987.2
$ws.Range("L40").Value = 3954987.2
$ws.Range("N40").Value = -3955337.2
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H49").Value = 475.7143
$ws.Range("I49").Value = 255
$ws.Range("J49").Value = 770
$ws.Range("K49").Value = 765
$ws.Range("L49").Value = 2310
$ws.Range("M49").Value = -629
$ws.Range("N49").Value = -2582
$ws.Range("H70").Value = 1071.4286
$ws.Range("J70").Value = 1500
$ws.Range("L70").Value = 4500
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 1071.4286
$ws.Range("J73").Value = 1500
$ws.Range("L73").Value = 4500
$ws.Range("N73").Value = -6372
$ws.Range("H80").Value = 6962666
$ws.Range("J80").Value = 11603969
$ws.Range("L80").Value = 34811907
$ws.Range("N80").Value = -34813903
$ws.Range("H82").Value = 327.8
$ws.Range("I82").Value = 327.8
$ws.Range("K82").Value = 983.4000000000001
$ws.Range("M82").Value = -577.4000000000001
$ws.Range("H83").Value = 6962666
$ws.Range("J83").Value = 11603969
$ws.Range("L83").Value = 104435721
$ws.Range("N83").Value = -104445705
$ws.Range("H85").Value = 327.8
$ws.Range("I85").Value = 327.8
$ws.Range("K85").Value = 983.4000000000001
$ws.Range("M85").Value = 420.5999999999999
$ws.Range("H86").Value = 18427.445
$ws.Range("I86").Value = 10601.2
$ws.Range("J86").Value = 28210.25
$ws.Range("K86").Value = 10601.2
$ws.Range("L86").Value = 28210.25
$ws.Range("M86").Value = -9478.200000000001
$ws.Range("N86").Value = -30456.25
$ws.Range("H89").Value = 18427.445
$ws.Range("I89").Value = 10601.2
$ws.Range("J89").Value = 28210.25
$ws.Range("K89").Value = 53006
$ws.Range("L89").Value = 141051.25
$ws.Range("M89").Value = -47390
$ws.Range("N89").Value = -152283.25
$ws.Range("H138").Value = 2391.0605
$ws.Range("J138").Value = 2289.1428
$ws.Range("L138").Value = 6867.428400000001
$ws.Range("N138").Value = -17147.4284

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 43479172
$ws.Range("I74").Value = 52632470
$ws.Range("J74").Value = 999.75
$ws.Range("K74").Value = 52632470
$ws.Range("L74").Value = 999.75
$ws.Range("M74").Value = -52631596
$ws.Range("N74").Value = -2747.75
$ws.Range("H77").Value = 43479172
$ws.Range("I77").Value = 52632470
$ws.Range("J77").Value = 999.75
$ws.Range("K77").Value = 263162350
$ws.Range("L77").Value = 4998.75
$ws.Range("M77").Value = -263157982
$ws.Range("N77").Value = -13734.75
$ws.Range("H88").Value = 143933
$ws.Range("I88").Value = 1375
$ws.Range("J88").Value = 334010.34
$ws.Range("K88").Value = 1375
$ws.Range("L88").Value = 334010.34
$ws.Range("M88").Value = -969
$ws.Range("N88").Value = -334822.34
$ws.Range("H91").Value = 143933
$ws.Range("I91").Value = 1375
$ws.Range("J91").Value = 334010.34
$ws.Range("K91").Value = 1375
$ws.Range("L91").Value = 334010.34
$ws.Range("M91").Value = 29
$ws.Range("N91").Value = -336818.34
$ws.Range("H132").Value = 28659
$ws.Range("J132").Value = 152918.9
$ws.Range("L132").Value = 458756.7
$ws.Range("N132").Value = -463816.7

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 909531.6
$ws.Range("I22").Value = 1250318.8
$ws.Range("K22").Value = 1250318.8
$ws.Range("M22").Value = -1250145.8
$ws.Range("H86").Value = 1355.4043
$ws.Range("I86").Value = 1245.2903
$ws.Range("J86").Value = 1568.75
$ws.Range("K86").Value = 1245.2903
$ws.Range("L86").Value = 1568.75
$ws.Range("M86").Value = -122.2902999999999
$ws.Range("N86").Value = -3814.75
$ws.Range("H89").Value = 1355.4043
$ws.Range("I89").Value = 1245.2903
$ws.Range("J89").Value = 1568.75
$ws.Range("K89").Value = 6226.451499999999
$ws.Range("L89").Value = 7843.75
$ws.Range("M89").Value = -610.4514999999992
$ws.Range("N89").Value = -19075.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 81714290
$ws.Range("I6").Value = 5500008
$ws.Range("J6").Value = 183333330
$ws.Range("K6").Value = 5500008
$ws.Range("L6").Value = 183333330
$ws.Range("M6").Value = -5499895
$ws.Range("N6").Value = -183333556
$ws.Range("H7").Value = 228
$ws.Range("I7").Value = 233.6
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 233.6
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -120.6
$ws.Range("N7").Value = -426
$ws.Range("H31").Value = 13050.389
$ws.Range("I31").Value = 14806.896
$ws.Range("K31").Value = 14806.896
$ws.Range("M31").Value = -14511.896
$ws.Range("H34").Value = 13050.389
$ws.Range("I34").Value = 14806.896
$ws.Range("K34").Value = 14806.896
$ws.Range("M34").Value = -14604.896
$ws.Range("H50").Value = 15990
$ws.Range("J50").Value = 15990
$ws.Range("L50").Value = 15990
$ws.Range("N50").Value = -17240
$ws.Range("H51").Value = 25600
$ws.Range("J51").Value = 25600
$ws.Range("L51").Value = 25600
$ws.Range("N51").Value = -27072
$ws.Range("H52").Value = 39899.5
$ws.Range("J52").Value = 39899.5
$ws.Range("L52").Value = 39899.5
$ws.Range("N52").Value = -40487.5
$ws.Range("H58").Value = 26444.4
$ws.Range("I58").Value = 1429.625
$ws.Range("J58").Value = 126503.5
$ws.Range("K58").Value = 1429.625
$ws.Range("L58").Value = 126503.5
$ws.Range("M58").Value = -1226.625
$ws.Range("N58").Value = -126909.5
$ws.Range("H61").Value = 25600
$ws.Range("J61").Value = 25600
$ws.Range("L61").Value = 25600
$ws.Range("N61").Value = -26296
$ws.Range("H132").Value = 12593.468
$ws.Range("I132").Value = 13506.561
$ws.Range("K132").Value = 40519.683
$ws.Range("M132").Value = -37989.683
$ws.Range("H136").Value = 26444.4
$ws.Range("I136").Value = 1429.625
$ws.Range("J136").Value = 126503.5
$ws.Range("K136").Value = 4288.875
$ws.Range("L136").Value = 379510.5
$ws.Range("M136").Value = -1738.875
$ws.Range("N136").Value = -384610.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 208.09091
$ws.Range("I18").Value = 148.77777
$ws.Range("J18").Value = 475
$ws.Range("K18").Value = 446.33331
$ws.Range("L18").Value = 1425
$ws.Range("M18").Value = -277.33331
$ws.Range("N18").Value = -1763
$ws.Range("H22").Value = 9818.182000000001
$ws.Range("I22").Value = 50400
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 151200
$ws.Range("L22").Value = 2400
$ws.Range("M22").Value = -151031
$ws.Range("N22").Value = -2738
$ws.Range("H27").Value = 9818.182000000001
$ws.Range("I27").Value = 50400
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 151200
$ws.Range("L27").Value = 2400
$ws.Range("M27").Value = -151098
$ws.Range("N27").Value = -2604
$ws.Range("H38").Value = 38461604
$ws.Range("I38").Value = 72.7
$ws.Range("K38").Value = 218.1
$ws.Range("M38").Value = 128.9
$ws.Range("H114").Value = 2265.4
$ws.Range("J114").Value = 2540
$ws.Range("L114").Value = 7620
$ws.Range("N114").Value = -14128
$ws.Range("H131").Value = 767.92
$ws.Range("J131").Value = 795.34045
$ws.Range("L131").Value = 2386.02135
$ws.Range("N131").Value = -12466.02135
$ws.Range("H132").Value = 1098.6666
$ws.Range("J132").Value = 1276.8572
$ws.Range("L132").Value = 11491.7148
$ws.Range("N132").Value = -16551.7148

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4160.6
$ws.Range("J80").Value = 4315.143
$ws.Range("L80").Value = 4315.143
$ws.Range("N80").Value = -6311.143
$ws.Range("H83").Value = 4160.6
$ws.Range("J83").Value = 4315.143
$ws.Range("L83").Value = 21575.715
$ws.Range("N83").Value = -31559.715
$ws.Range("H132").Value = 54936.07
$ws.Range("I132").Value = 53867.95
$ws.Range("J132").Value = 57309.668
$ws.Range("K132").Value = 161603.85
$ws.Range("L132").Value = 171929.004
$ws.Range("M132").Value = -159073.85
$ws.Range("N132").Value = -176989.004

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5296.8335
$ws.Range("I22").Value = 5220.25
$ws.Range("J22").Value = 5450
$ws.Range("K22").Value = 5220.25
$ws.Range("L22").Value = 5450
$ws.Range("M22").Value = -4925.25
$ws.Range("N22").Value = -6040
$ws.Range("H27").Value = 5296.8335
$ws.Range("I27").Value = 5220.25
$ws.Range("J27").Value = 5450
$ws.Range("K27").Value = 5220.25
$ws.Range("L27").Value = 5450
$ws.Range("M27").Value = -5113.25
$ws.Range("N27").Value = -5664
$ws.Range("H82").Value = 2240.5908
$ws.Range("I82").Value = 2105.8823
$ws.Range("J82").Value = 2698.6
$ws.Range("K82").Value = 2105.8823
$ws.Range("L82").Value = 2698.6
$ws.Range("M82").Value = -1744.8823
$ws.Range("N82").Value = -3420.6
$ws.Range("H85").Value = 2240.5908
$ws.Range("I85").Value = 2105.8823
$ws.Range("J85").Value = 2698.6
$ws.Range("K85").Value = 2105.8823
$ws.Range("L85").Value = 2698.6
$ws.Range("M85").Value = -857.8823000000002
$ws.Range("N85").Value = -5194.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58824316
$ws.Range("I81").Value = 777.3333
$ws.Range("J81").Value = 125000800
$ws.Range("K81").Value = 1554.6666
$ws.Range("L81").Value = 250001600
$ws.Range("M81").Value = -493.6666
$ws.Range("N81").Value = -250003722
$ws.Range("H84").Value = 58824316
$ws.Range("I84").Value = 777.3333
$ws.Range("J84").Value = 125000800
$ws.Range("K84").Value = 7773.333000000001
$ws.Range("L84").Value = 1250008000
$ws.Range("M84").Value = -2469.333000000001
$ws.Range("N84").Value = -1250018608
$ws.Range("H132").Value = 1379.8096
$ws.Range("I132").Value = 1004.8823
$ws.Range("K132").Value = 3014.6469
$ws.Range("M132").Value = -484.6468999999997

Write-Output "Applied 259 cell edits across 8 sheets."